# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 13:33"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1213010
$ws.Range("C4").Value = 175
$ws.Range("E4").Value = 955017

# --- Row 62: Barein ---
$ws.Range("B62").Value = 3679
$ws.Range("C62").Value = 146
$ws.Range("D62").Value = 1762
$ws.Range("E62").Value = 1909

# --- Rows 112/113: San Marino moves above El Salvador in the sorted list ---
# Row 112 becomes San Marino (updated data), row 113 becomes El Salvador
# (its data reverts to what was previously El Salvador's row).
$ws.Range("A112").Value = "San Marino"
$ws.Range("B112").Value = 589
$ws.Range("C112").Value = 7
$ws.Range("D112").Value = 92
$ws.Range("E112").Value = 456
$ws.Range("F112").Value = 5
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 41

$ws.Range("A113").Value = "El Salvador"
$ws.Range("B113").Value = 587
$ws.Range("C113").Value = 32
$ws.Range("D113").Value = 201
$ws.Range("E113").Value = 373
$ws.Range("F113").Value = 3
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 13
